$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, pushing all existing game rows down by one.
$ws.Rows.Item(1).Insert()

# Populate the new header row with the column titles.
$ws.Cells.Item(1, 1).Value = "Title"
$ws.Cells.Item(1, 2).Value = "URL"
$ws.Cells.Item(1, 3).Value = "Status"
